$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top; this shifts all existing rows
# (including the old header row 1) down by two.
$ws.Rows("1:2").Insert()

# Copy the formatting (style) of the former header row -- now row 3 --
# onto the new row 1 so it keeps the bold/bordered/centered header style.
$ws.Range("A3:N3").Copy()
$ws.Range("A1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 1: a simple numeric index sequence 0..13 across columns A..N.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10
$ws.Range("L1").Value = 11
$ws.Range("M1").Value = 12
$ws.Range("N1").Value = 13

# New row 2: blank across the board, except column E which reads "Drive".
$ws.Range("E2").Value = "Drive"
